$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.191.64"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").Value = "3.727.71"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.10%  "

$ws.Range("E7").Value = "  +1.08%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("E9").Value = "  +2.31%  "

$ws.Range("E10").Value = "  -2.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000291"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("D14").Value = "4.330.08"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").Value = "3.730.83"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("D20").Value = "69.085.15"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("E22").Value = "  +1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "90.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.39%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "

$ws.Range("E33").Value = "  +4.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "641.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "45.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.85%  "

$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.417"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.66%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0834"
$ws.Range("E38").Value = "  -8.82%  "

$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("E41").Value = "  +3.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0450"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.13%  "

$ws.Range("D45").Value = "2.906.15"
$ws.Range("E45").Value = "  +5.26%  "

$ws.Range("E46").Value = "  +3.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.25%  "

$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("E51").Value = "  -13.59%  "
